$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Make room for three new rows of metadata above the existing content
# (national_indicator_description / other_info / data_show_map).
# ------------------------------------------------------------------
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# ------------------------------------------------------------------
# Populate the three new rows.
# ------------------------------------------------------------------
$ws.Range("A1").Value2 = "SDG_GOAL"
$ws.Range("B1").Value2 = "Bevölkerung"

$ws.Range("A2").Value2 = "indicator_name"
$ws.Range("B2").Value2 = "Bevölkerungsentwicklung in Niedersachsen und den Bundesländern"

$ws.Range("A3").Value2 = "SDG_INDICATOR"
$ws.Range("B3").Value2 = "Bev. Entwicklung"

# ------------------------------------------------------------------
# Give A1:A3 the same "key" look as the other label cells in column A
# (7pt Consolas) by copying the format from the data_show_map label,
# then recolor it so a brand-new font/style entry is produced.
# ------------------------------------------------------------------
$ws.Range("A6").Copy()
$ws.Range("A1:A3").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Font.Color = 2712337
$ws.Range("A1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# Restore the tall wrapped-text row height for the long description
# cell, which is now row 4 after the insert.
# ------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 250.65

# ------------------------------------------------------------------
# Widen column B so the new long text values are readable.
# (ColumnWidth is internally snapped to the nearest 1/6, so this lands
# on the closest achievable stored width to 104.5546875.)
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 103.65

# ------------------------------------------------------------------
# Reset the view: select A3 instead of the old B3 selection, and drop
# the scrolled-down topLeftCell.
# ------------------------------------------------------------------
$ws.Range("A3").Select() | Out-Null

Write-Output "done"
